# "Generate Report for Handoff"
#
# The handoff transform failed for the 7b2d39ff-...md source file, so:
#   - the "Ready for handoff" status everywhere it appears becomes
#     "Handoff transform failed" (Overview!B2/C2, zh-cn!B2, de-de!B2)
#   - the per-language detail sheets (zh-cn, de-de) lose their "Latest
#     Handoff File" link/value (column C) and their handoff datetime /
#     handback-datetime / handoff-reason columns reset to the "nothing
#     happened yet" defaults (0001-01-01 00:00:00 / Ignored)

$wb = $excel.ActiveWorkbook

# --- Overview sheet: the shared "Ready for handoff" status text changes ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = "Handoff transform failed"
$overview.Range("C2").Value = "Handoff transform failed"

# --- Per-language detail sheets (zh-cn / de-de) ---
$langSheets = @(
    @{ Name = "zh-cn"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4fcb51403be1da82d473ec9db602429482408b1b/e2e/7b2d39ff-d269-4e15-87bc-9b12eafbb1f9.md"; ConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4fcb51403be1da82d473ec9db602429482408b1b/.localization-config" },
    @{ Name = "de-de"; MdUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4fcb51403be1da82d473ec9db602429482408b1b/e2e/7b2d39ff-d269-4e15-87bc-9b12eafbb1f9.md"; ConfigUrl = "https://github.com/OpenLocalizationTest/oltest/blob/4fcb51403be1da82d473ec9db602429482408b1b/.localization-config" }
)

foreach ($entry in $langSheets) {
    $ws = $wb.Worksheets.Item($entry.Name)

    # Status text: "Ready for handoff" -> "Handoff transform failed"
    $ws.Range("B2").Value = "Handoff transform failed"

    # The "Latest Handoff File" cell (C2) and its hyperlink go away entirely
    # - the transform never produced a handoff file this run.
    $ws.Range("C2").ClearContents()

    # Latest Handoff Datetime / Handback DateTime / Handoff Reason reset
    $ws.Range("D2").Value = "0001-01-01 00:00:00"
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Ignored"

    $ws.Range("D3").Value = "0001-01-01 00:00:00"
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Ignored"

    # Hyperlinks.Delete() on this engine drops every hyperlink on the sheet,
    # so rebuild the two that must survive (A2 -> the .md source file, A3 ->
    # the .localization-config file) after clearing them all; this also
    # drops the now-orphaned C2 hyperlink and renumbers the relationship ids
    # the same way the target workbook does (rId3 free'd up by removing the
    # C2 link becomes A3's new hyperlink id).
    $ws.Hyperlinks.Delete()
    $ws.Hyperlinks.Add($ws.Range("A2"), $entry.MdUrl, "", "", "7b2d39ff-d269-4e15-87bc-9b12eafbb1f9.md")
    $ws.Hyperlinks.Add($ws.Range("A3"), $entry.ConfigUrl, "", "", ".localization-config")
}
